$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics (per commit "update scripts wuth new tpm")

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.609586333333333
$ws.Range("H2").Value2 = 4.828759
$ws.Range("I2").Value2 = 0.05107819292772156
$ws.Range("J2").Value2 = 0.05107819292772156
$ws.Range("M2").Value2 = 201.098592
$ws.Range("N2").Value2 = 603.295776
$ws.Range("O2").Value2 = 0.7918622805845071
$ws.Range("P2").Value2 = 0.791862280584507
$ws.Range("Q2").Value2 = 323.6855453357761
$ws.Range("R2").Value2 = 2913.169908021984
$ws.Range("S2").Value2 = 0.04044689433988104
$ws.Range("T2").Value2 = 0.04044689433988103

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.609586333333333
$ws.Range("H3").Value2 = 4.828759
$ws.Range("I3").Value2 = 0.05107819292772156
$ws.Range("J3").Value2 = 0.05107819292772156
$ws.Range("O3").Value2 = 0.1414593902976603
$ws.Range("P3").Value2 = 0.1414593902976603
$ws.Range("Q3").Value2 = 57.823640567355
$ws.Range("R3").Value2 = 520.412765106195
$ws.Range("S3").Value2 = 0.007225490029061755
$ws.Range("T3").Value2 = 0.007225490029061755

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.609586333333333
$ws.Range("H4").Value2 = 4.828759
$ws.Range("I4").Value2 = 0.05107819292772156
$ws.Range("J4").Value2 = 0.05107819292772156
$ws.Range("O4").Value2 = 0.0666783291178327
$ws.Range("P4").Value2 = 0.06667832911783268
$ws.Range("Q4").Value2 = 27.25576385157889
$ws.Range("R4").Value2 = 245.30187466421
$ws.Range("S4").Value2 = 0.003405808558778773
$ws.Range("T4").Value2 = 0.003405808558778772

# Row 5
$ws.Range("I5").Value2 = 0.5992082897496871
$ws.Range("J5").Value2 = 0.5992082897496871
$ws.Range("M5").Value2 = 201.098592
$ws.Range("N5").Value2 = 603.295776
$ws.Range("O5").Value2 = 0.7918622805845071
$ws.Range("P5").Value2 = 0.791862280584507
$ws.Range("Q5").Value2 = 3797.218556885953
$ws.Range("R5").Value2 = 34174.96701197357
$ws.Range("S5").Value2 = 0.4744904428663294
$ws.Range("T5").Value2 = 0.4744904428663293

# Row 6
$ws.Range("I6").Value2 = 0.5992082897496871
$ws.Range("J6").Value2 = 0.5992082897496871
$ws.Range("O6").Value2 = 0.1414593902976603
$ws.Range("P6").Value2 = 0.1414593902976603
$ws.Range("S6").Value2 = 0.0847636393292945
$ws.Range("T6").Value2 = 0.0847636393292945

# Row 7
$ws.Range("I7").Value2 = 0.5992082897496871
$ws.Range("J7").Value2 = 0.5992082897496871
$ws.Range("O7").Value2 = 0.0666783291178327
$ws.Range("P7").Value2 = 0.06667832911783268
$ws.Range("S7").Value2 = 0.03995420755406329
$ws.Range("T7").Value2 = 0.03995420755406329

# Row 8
$ws.Range("H8").Value2 = 33.060729
$ws.Range("I8").Value2 = 0.3497135173225914
$ws.Range("J8").Value2 = 0.3497135173225914
$ws.Range("M8").Value2 = 201.098592
$ws.Range("N8").Value2 = 603.295776
$ws.Range("O8").Value2 = 0.7918622805845071
$ws.Range("P8").Value2 = 0.791862280584507
$ws.Range("Q8").Value2 = 2216.155350797856
$ws.Range("R8").Value2 = 19945.39815718071
$ws.Range("S8").Value2 = 0.2769249433782968
$ws.Range("T8").Value2 = 0.2769249433782968

# Row 9
$ws.Range("H9").Value2 = 33.060729
$ws.Range("I9").Value2 = 0.3497135173225914
$ws.Range("J9").Value2 = 0.3497135173225914
$ws.Range("O9").Value2 = 0.1414593902976603
$ws.Range("P9").Value2 = 0.1414593902976603
$ws.Range("Q9").Value2 = 395.897105362005
$ws.Range("S9").Value2 = 0.04947026093930404
$ws.Range("T9").Value2 = 0.04947026093930404

# Row 10
$ws.Range("H10").Value2 = 33.060729
$ws.Range("I10").Value2 = 0.3497135173225914
$ws.Range("J10").Value2 = 0.3497135173225914
$ws.Range("O10").Value2 = 0.0666783291178327
$ws.Range("P10").Value2 = 0.06667832911783268
$ws.Range("Q10").Value2 = 186.6101460820567
$ws.Range("S10").Value2 = 0.02331831300499064
$ws.Range("T10").Value2 = 0.02331831300499063
